# "Fruta / hortaliza, semanal"
# Insert a new weekly data row for Zanahoria (Terminal Hortofrutícola Agro
# Chillán) at row 289, pushing the existing rows 289-349 down to 290-350.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 289 (shifts 289..349 -> 290..350)
$ws.Rows.Item(289).Insert()

# Populate the new row 289 with the latest weekly price observation
$ws.Range("A289").Value = 7
$ws.Range("B289").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C289").Value = "Ñuble"
$ws.Range("D289").Value = 44889
$ws.Range("E289").Value = 16
$ws.Range("F289").Value = 100114013
$ws.Range("G289").Value = "Zanahoria"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 160
$ws.Range("K289").Value = 9000
$ws.Range("L289").Value = 10000
$ws.Range("M289").Value = 9500
$ws.Range("N289").Value = "$/saco 20 kilos"
$ws.Range("O289").Value = "Región de Ñuble"
$ws.Range("P289").Value = 475
$ws.Range("Q289").Value = 20
$ws.Range("R289").Value = "Hortaliza"
